# Update the "Jquery ... scripting" bullet on the "Our Design Choices" slide.
# Before: "Jquery" + " for client-side scripting"
# After:  "Jquery" + " as a client-side " + "scripting library"   (as two runs)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$oldTail = " for client-side scripting"
$newFirst = " as a client-side "
$newSecond = "scripting library"

# Locate the run to edit, dynamically, so the script does not depend on
# hard-coded character offsets.
$full = $tr.Text
$idx = $full.IndexOf($oldTail)
if ($idx -lt 0) {
    throw "Could not find target text '$oldTail' in shape."
}
$start = $idx + 1
$len = $oldTail.Length

# Replace the whole old run's text with the full new text first (keeps the
# original run's formatting / dirty="0" on the first part).
$target = $tr.Characters($start, $len)
$target.Text = $newFirst + $newSecond

# Re-find the "scripting library" portion and re-assign its text so that it
# is split off into its own run (second run of the pair).
$trAfter = $sh.TextFrame.TextRange
$fullAfter = $trAfter.Text
$idx2 = $fullAfter.IndexOf($newSecond, $idx)
$start2 = $idx2 + 1
$len2 = $newSecond.Length
$secondRun = $trAfter.Characters($start2, $len2)
$secondRun.Text = $newSecond
